# Applies the "Added more to livrabil" edit:
#  1. Expands the placeholder "[ Descriere despre joc ]" into the real
#     project-scope paragraph (split across several same-format runs,
#     exactly as the author's Word session produced them).
#  2. Fixes the "intefata" typo to "interfata" inside the bold run
#     "intefata grafica" (splitting it into three bold runs).
#  3. Moves the <w:lastRenderedPageBreak/> marker: it now lands on the
#     "O functionalitate ..." paragraph instead of the
#     "Va exista si o componenta de ..." paragraph (pagination shifted
#     because of the extra text added above).
#  4. Likewise moves a second <w:lastRenderedPageBreak/> marker from the
#     "Un user poate crea o legatura ..." paragraph onto the
#     "Userii sunt initial off-line ..." paragraph.

$d = $word.ActiveDocument

$wordMl = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($marker, $innerXml) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$marker*") {
            $rng = $p.Range
            $rng.Text = ""
            $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wordMl + '><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $rng.InsertXML($xml)
            return
        }
    }
    throw "paragraph containing '$marker' not found"
}

# 1) "[ Descriere despre joc ]" -> full scope description
$p1 = @'
<w:p w14:paraId="0243E537" w14:textId="4414C319" w:rsidR="00241B79" w:rsidRPr="00241B79" w:rsidRDefault="00650C37" w:rsidP="00241B79"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">Programul are ca scop realizarea unui joc de tip Bomberman </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>î</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>ntr-un mod intuitiv</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> usor de utilizat. Accept</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>nd mai multi clienti, fiecare client isi poate personaliza experienta, alegand din lista de jucatori online, avand in acelasi timp access si la un istoric al activitatii</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>
'@
Set-ParagraphXml "Descriere despre joc" $p1

# 2) "intefata grafica" -> "interfata grafica" (typo fix, still bold, split
#    into three runs: "inte" / "r" / "fata grafica")
$p2 = @'
<w:p w14:paraId="7BA368DF" w14:textId="436E6AA4" w:rsidR="00241B79" w:rsidRDefault="00241B79" w:rsidP="00241B79"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">Clientii vor comunica cu server-ul prin intermediul unei aplicatii pentru </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>desktop</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">, ce dispune de o </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>inte</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>fata grafica</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> in care va fi afisat jocul.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Clientii vor comunica" $p2

# 3) "O functionalitate ..." paragraph gains the lastRenderedPageBreak
$p3 = @'
<w:p w14:paraId="7DFDF5B5" w14:textId="5173918D" w:rsidR="00241B79" w:rsidRDefault="00241B79" w:rsidP="00241B79"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:lastRenderedPageBreak/><w:t>O functionalitate pe care o mai are server-ul este ca, la inchidere isi va salva in mod automat datele in cadrul unei baze de date, de exemplu pentru utilizatori</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>id unic, tip utilizator, username</w:t></w:r><w:r w:rsidR="00BA4A98"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-ParagraphXml "O functionalitate pe care" $p3

# 4) "Va exista si o componenta de ..." paragraph loses its lastRenderedPageBreak
$p4 = @'
<w:p w14:paraId="59A39F72" w14:textId="1027D3A4" w:rsidR="00BA4A98" w:rsidRDefault="00BA4A98" w:rsidP="00BA4A98"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">Va exista si o componenta de </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>logging</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">, care va marca intr-un fisier fiecare actiune pe care o va executa server-ul. Numele fisierului va fi precizat la initializare server-ului prin intermediul unui </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>fisier de configurare</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Va exista si o componenta de" $p4

# 5) "Userii sunt initial off-line ..." paragraph gains the lastRenderedPageBreak
$p5 = @'
<w:p w14:paraId="4BDF9712" w14:textId="749964B8" w:rsidR="001844C9" w:rsidRDefault="001844C9" w:rsidP="00BA4A98"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:lastRenderedPageBreak/><w:t>Userii sunt initial off-line. Pentru a beneficia de serviciile sistemului,un user trebuie sa devina online, realizand autentificarea la server.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Userii sunt initial off-line" $p5

# 6) "Un user poate crea o legatura ..." paragraph loses its lastRenderedPageBreak
$p6 = @'
<w:p w14:paraId="174F877B" w14:textId="12C8CB38" w:rsidR="001844C9" w:rsidRDefault="001844C9" w:rsidP="00BA4A98"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Un user poate crea o legatura cu un alt user prin intermediul username-ului.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Un user poate crea o legatura" $p6

Write-Output "edit applied"
